# Finish experiments for the day
# Fill in the results for the last three Instrument5 runs (rows 76-78 on the
# "Log" sheet) and record the logfile/date for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Row 76 already has its Accuracy (C76) recorded; just log the date/file.
$ws.Range("G76").Value = 42557
$ws.Range("H76").Value = "16.29.txt"

# Row 77: record the accuracy plus date/logfile.
$ws.Range("C77").Value = 0.344
$ws.Range("G77").Value = 42557
$ws.Range("H77").Value = "16.29.txt"

# Row 78: record the accuracy plus date/logfile.
$ws.Range("C78").Value = 0.368
$ws.Range("G78").Value = 42557
$ws.Range("H78").Value = "16.29.txt"

# Move the active selection to the last cell touched, like a user would
# after finishing data entry for the day.
$ws.Range("H78").Select()
